# Update shared-string content of cell E8 on the "Rules" sheet from
# "Good Morning" to "GIT UPDATE", and leave that cell selected/active
# (matches the commit "update file with jgit").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Activate()

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
